# "Implemented auto-update shop UI"
# Adds the latest hour-registration entry (row 26) plus the mirrored
# "shop" tracking entry in the right-hand table (row 4, columns O-S).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Right-hand table: new "shop" activity entry on row 4 ---
$ws.Range("O4").Value = 44722
$ws.Range("O4").NumberFormat = "d-mmm"
$ws.Range("P4").Value = "12.00 - 13.00"
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = "Tested object pooling, tried out occlusion culling"

# --- Main hour-registration table: new entry, row 26 ---
$ws.Range("A26").Value = 44722
$ws.Range("A26").NumberFormat = "d-mmm"
$ws.Range("B26").Value = "10.30 - 12.00"
$ws.Range("D26").Value = 1.5
$ws.Range("E26").Value = "Polishing gunplay"

# Move the selection to the newly added row, matching the author's last
# edit position.
$ws.Range("F26").Select() | Out-Null
